$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed new shared strings in the exact order they were introduced, so the
# --- resulting shared-strings table allocation matches the target workbook. ---
$ws.Range("P47").Value = "Compare whiten and non whiten"
$ws.Range("Q49").Value = "RW BS15_AU_02a_files_1-104"
$ws.Range("Q50").Value = "RW BS14_AU_04_files_All"
$ws.Range("Q51").Value = "RW AW12_AU_BS3_files_All"
$ws.Range("Q52").Value = "RW BS13_AU_04_files_All"
$ws.Range("Q53").Value = "RW BS16_AU_02a_files_1-175"
$ws.Range("Q54").Value = "RW BS15_AU_02b_files_All"
$ws.Range("Q55").Value = "RW AW14_AU_BS3_files_1-160"
$ws.Range("P57").Value = "all"
$ws.Range("P48").Value = "non whiten: "
$ws.Range("Y56").Value = "* after raven GT'd "
$ws.Range("P58").Value = "whiten"
$ws.Range("Y62").Value = "arguably did worse "
$ws.Range("O49").Value = 4
$ws.Range("P49").Value = "RW"
$ws.Range("R49").Value = 1559
$ws.Range("S49").Value = 276
$ws.Range("T49").Value = 1283
$ws.Range("U49").Value = 6
$ws.Range("V49").Value = 0.97872340425531901
$ws.Range("W49").Value = 0.82296343810134698
$ws.Range("X49").Value = 0.215120810600156
$ws.Range("O50").Value = 4
$ws.Range("P50").Value = "RW"
$ws.Range("R50").Value = 5015
$ws.Range("S50").Value = 434
$ws.Range("T50").Value = 4581
$ws.Range("U50").Value = 29
$ws.Range("V50").Value = 0.93736501079913603
$ws.Range("W50").Value = 0.91345962113658996
$ws.Range("X50").Value = 0.094739139925780397
$ws.Range("O51").Value = 4
$ws.Range("P51").Value = "RW"
$ws.Range("R51").Value = 5346
$ws.Range("S51").Value = 805
$ws.Range("T51").Value = 4541
$ws.Range("U51").Value = 50
$ws.Range("V51").Value = 0.94152046783625698
$ws.Range("W51").Value = 0.84942012719790505
$ws.Range("X51").Value = 0.17727372825368901
$ws.Range("O52").Value = 4
$ws.Range("P52").Value = "RW"
$ws.Range("R52").Value = 12191
$ws.Range("S52").Value = 1477
$ws.Range("T52").Value = 10714
$ws.Range("U52").Value = 129
$ws.Range("V52").Value = 0.91967621419676204
$ws.Range("W52").Value = 0.87884504962677401
$ws.Range("X52").Value = 0.13785700952025401
$ws.Range("O53").Value = 4
$ws.Range("P53").Value = "RW"
$ws.Range("R53").Value = 3174
$ws.Range("S53").Value = 471
$ws.Range("T53").Value = 2703
$ws.Range("U53").Value = 19
$ws.Range("V53").Value = 0.96122448979591801
$ws.Range("W53").Value = 0.85160680529300603
$ws.Range("X53").Value = 0.174250832408435
$ws.Range("O54").Value = 4
$ws.Range("P54").Value = "RW"
$ws.Range("R54").Value = 2528
$ws.Range("S54").Value = 160
$ws.Range("T54").Value = 2368
$ws.Range("U54").Value = 19
$ws.Range("V54").Value = 0.89385474860335201
$ws.Range("W54").Value = 0.936708860759494
$ws.Range("X54").Value = 0.067567567567567599
$ws.Range("O55").Value = 4
$ws.Range("P55").Value = "RW"
$ws.Range("R55").Value = 5641
$ws.Range("S55").Value = 570
$ws.Range("T55").Value = 5071
$ws.Range("U55").Value = 41
$ws.Range("V55").Value = 0.93289689034369905
$ws.Range("W55").Value = 0.89895408615493699
$ws.Range("X55").Value = 0.112403865115362
$ws.Range("O56").Value = 4
$ws.Range("P56").Value = "RW"
$ws.Range("Q56").Value = "RW AL16_AU_BS1_files_All"
$ws.Range("R56").Value = 5880
$ws.Range("S56").Value = 440
$ws.Range("T56").Value = 5440
$ws.Range("U56").Value = 0
$ws.Range("V56").Value = 1
$ws.Range("W56").Value = 0.92517006802721102
$ws.Range("X56").Value = 0.080882352941176502
$ws.Range("O57").Value = 4
$ws.Range("R57").Value = 41334
$ws.Range("S57").Value = 4633
$ws.Range("T57").Value = 36701
$ws.Range("U57").Value = 293
$ws.Range("V57").Value = 0.94051969099999999
$ws.Range("W57").Value = 0.88791309799999996
$ws.Range("X57").Value = 0.126236342
$ws.Range("O59").Value = 4
$ws.Range("P59").Value = "RW"
$ws.Range("Q59").Value = "RW BS15_AU_02a_files_1-104"
$ws.Range("R59").Value = 1496
$ws.Range("S59").Value = 278
$ws.Range("T59").Value = 1218
$ws.Range("U59").Value = 3
$ws.Range("V59").Value = 0.98932384341637003
$ws.Range("W59").Value = 0.814171122994652
$ws.Range("X59").Value = 0.22824302134647001
$ws.Range("O60").Value = 4
$ws.Range("P60").Value = "RW"
$ws.Range("Q60").Value = "RW BS14_AU_04_files_All"
$ws.Range("R60").Value = 4905
$ws.Range("S60").Value = 441
$ws.Range("T60").Value = 4464
$ws.Range("U60").Value = 23
$ws.Range("V60").Value = 0.95043103448275901
$ws.Range("W60").Value = 0.91009174311926599
$ws.Range("X60").Value = 0.098790322580645198
$ws.Range("O61").Value = 4
$ws.Range("P61").Value = "RW"
$ws.Range("Q61").Value = "RW AW12_AU_BS3_files_All"
$ws.Range("R61").Value = 5346
$ws.Range("S61").Value = 810
$ws.Range("T61").Value = 4536
$ws.Range("U61").Value = 46
$ws.Range("V61").Value = 0.94626168224299101
$ws.Range("W61").Value = 0.84848484848484895
$ws.Range("X61").Value = 0.17857142857142899
$ws.Range("O62").Value = 4
$ws.Range("P62").Value = "RW"
$ws.Range("Q62").Value = "RW BS13_AU_04_files_All"
$ws.Range("R62").Value = 12297
$ws.Range("S62").Value = 1482
$ws.Range("T62").Value = 10815
$ws.Range("U62").Value = 128
$ws.Range("V62").Value = 0.92049689440993798
$ws.Range("W62").Value = 0.87948280068309304
$ws.Range("X62").Value = 0.13703190013869601
$ws.Range("O63").Value = 4
$ws.Range("P63").Value = "RW"
$ws.Range("Q63").Value = "RW BS16_AU_02a_files_1-175"
$ws.Range("R63").Value = 3168
$ws.Range("S63").Value = 470
$ws.Range("T63").Value = 2698
$ws.Range("U63").Value = 20
$ws.Range("V63").Value = 0.95918367346938804
$ws.Range("W63").Value = 0.85164141414141403
$ws.Range("X63").Value = 0.17420311341734601
$ws.Range("O64").Value = 4
$ws.Range("P64").Value = "RW"
$ws.Range("Q64").Value = "RW BS15_AU_02b_files_All"
$ws.Range("R64").Value = 2400
$ws.Range("S64").Value = 161
$ws.Range("T64").Value = 2239
$ws.Range("U64").Value = 17
$ws.Range("V64").Value = 0.90449438202247201
$ws.Range("W64").Value = 0.93291666666666695
$ws.Range("X64").Value = 0.071907101384546701
$ws.Range("O65").Value = 4
$ws.Range("P65").Value = "RW"
$ws.Range("Q65").Value = "RW AW14_AU_BS3_files_1-160"
$ws.Range("R65").Value = 5539
$ws.Range("S65").Value = 567
$ws.Range("T65").Value = 4972
$ws.Range("U65").Value = 44
$ws.Range("V65").Value = 0.92798690671031103
$ws.Range("W65").Value = 0.89763495215742894
$ws.Range("X65").Value = 0.114038616251006
$ws.Range("O66").Value = 4
$ws.Range("P66").Value = "RW"
$ws.Range("Q66").Value = "RW AL16_AU_BS1_files_All"
$ws.Range("R66").Value = 5947
$ws.Range("S66").Value = 430
$ws.Range("T66").Value = 5517
$ws.Range("U66").Value = 13
$ws.Range("V66").Value = 0.97065462753950305
$ws.Range("W66").Value = 0.92769463595089996
$ws.Range("X66").Value = 0.077940909914808801
$ws.Range("O67").Value = 4
$ws.Range("P67").Value = "RW"
$ws.Range("Q67").Value = "all"
$ws.Range("R67").Value = 41098
$ws.Range("S67").Value = 4639
$ws.Range("T67").Value = 36459
$ws.Range("U67").Value = 294
$ws.Range("V67").Value = 0.94040137847151795
$ws.Range("W67").Value = 0.88712346099566897
$ws.Range("X67").Value = 0.127238816204504
$ws.Range("O69").Value = 4
$ws.Range("P69").Value = "RW"
$ws.Range("Q69").Value = "rf BS15_AU_02a_files_1-104"
$ws.Range("R69").Value = 477
$ws.Range("S69").Value = 260
$ws.Range("T69").Value = 217
$ws.Range("U69").Value = 20
$ws.Range("V69").Value = 0.92857142857142905
$ws.Range("W69").Value = 0.45492662473794498
$ws.Range("X69").Value = 1.1981566820276499
$ws.Range("Y69").Value = 0.96553217458657603
$ws.Range("O70").Value = 4
$ws.Range("P70").Value = "RW"
$ws.Range("Q70").Value = "rf BS14_AU_04_files_All"
$ws.Range("R70").Value = 1229
$ws.Range("S70").Value = 390
$ws.Range("T70").Value = 839
$ws.Range("U70").Value = 71
$ws.Range("V70").Value = 0.84598698481561796
$ws.Range("W70").Value = 0.68266883645239995
$ws.Range("X70").Value = 0.46483909415971397
$ws.Range("Y70").Value = 0.96553217458657603
$ws.Range("O71").Value = 4
$ws.Range("P71").Value = "RW"
$ws.Range("Q71").Value = "rf AW12_AU_BS3_files_All"
$ws.Range("R71").Value = 1681
$ws.Range("S71").Value = 752
$ws.Range("T71").Value = 929
$ws.Range("U71").Value = 101
$ws.Range("V71").Value = 0.88159437280187603
$ws.Range("W71").Value = 0.55264723378941105
$ws.Range("X71").Value = 0.80947255113024796
$ws.Range("Y71").Value = 0.96553217458657603
$ws.Range("O72").Value = 4
$ws.Range("P72").Value = "RW"
$ws.Range("Q72").Value = "rf BS13_AU_04_files_All"
$ws.Range("R72").Value = 3294
$ws.Range("S72").Value = 1336
$ws.Range("T72").Value = 1958
$ws.Range("U72").Value = 268
$ws.Range("V72").Value = 0.83291770573566104
$ws.Range("W72").Value = 0.59441408621736502
$ws.Range("X72").Value = 0.68232890704800797
$ws.Range("Y72").Value = 0.96553217458657603
$ws.Range("O73").Value = 4
$ws.Range("P73").Value = "RW"
$ws.Range("Q73").Value = "rf BS16_AU_02a_files_1-175"
$ws.Range("R73").Value = 1020
$ws.Range("S73").Value = 441
$ws.Range("T73").Value = 579
$ws.Range("U73").Value = 49
$ws.Range("V73").Value = 0.9
$ws.Range("W73").Value = 0.56764705882352895
$ws.Range("X73").Value = 0.76165803108808305
$ws.Range("Y73").Value = 0.96553217458657603
$ws.Range("O74").Value = 4
$ws.Range("P74").Value = "RW"
$ws.Range("Q74").Value = "rf BS15_AU_02b_files_All"
$ws.Range("R74").Value = 579
$ws.Range("S74").Value = 143
$ws.Range("T74").Value = 436
$ws.Range("U74").Value = 34
$ws.Range("V74").Value = 0.80790960451977401
$ws.Range("W74").Value = 0.75302245250431799
$ws.Range("X74").Value = 0.32798165137614699
$ws.Range("Y74").Value = 0.96553217458657603
$ws.Range("O75").Value = 4
$ws.Range("P75").Value = "RW"
$ws.Range("Q75").Value = "rf AW14_AU_BS3_files_1-160"
$ws.Range("R75").Value = 1621
$ws.Range("S75").Value = 514
$ws.Range("T75").Value = 1107
$ws.Range("U75").Value = 92
$ws.Range("V75").Value = 0.84818481848184801
$ws.Range("W75").Value = 0.68291178285009302
$ws.Range("X75").Value = 0.46431797651309797
$ws.Range("Y75").Value = 0.96553217458657603
$ws.Range("O76").Value = 4
$ws.Range("P76").Value = "RW"
$ws.Range("Q76").Value = "rf AL16_AU_BS1_files_All"
$ws.Range("R76").Value = 1449
$ws.Range("S76").Value = 376
$ws.Range("T76").Value = 1073
$ws.Range("U76").Value = 63
$ws.Range("V76").Value = 0.85649202733485197
$ws.Range("W76").Value = 0.74051069703243599
$ws.Range("X76").Value = 0.35041938490214403
$ws.Range("Y76").Value = 0.96553217458657603
$ws.Range("O77").Value = 4
$ws.Range("P77").Value = "RW"
$ws.Range("Q77").Value = "rf all"
$ws.Range("R77").Value = 11350
$ws.Range("S77").Value = 4212
$ws.Range("T77").Value = 7138
$ws.Range("U77").Value = 698
$ws.Range("V77").Value = 0.85784114052953198
$ws.Range("W77").Value = 0.62889867841409697
$ws.Range("X77").Value = 0.59008125525357202
$ws.Range("Y77").Value = 0.96553217458657603

# --- Restore the scroll position / active selection recorded in the sheet view ---
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("U71").Select()
